# fix 9mm damage 2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- bullet_damage (I) tweaks ---
$ws.Range("I7").Value = 0.09
$ws.Range("I9").Value = 0.03
$ws.Range("I14").Value = 0.03
$ws.Range("I15").Value = 0.04
$ws.Range("I16").Value = 0.09
$ws.Range("I17").Value = 0.09
$ws.Range("I18").Value = 0.09
$ws.Range("I20").Value = 0.04
$ws.Range("I21").Value = 0.07
$ws.Range("I23").Value = 0.09
$ws.Range("I24").Value = 0.06
$ws.Range("I26").Value = 0.03
$ws.Range("I27").Value = 0.03
$ws.Range("I28").Value = 0.09
$ws.Range("I31").Value = 0.02

# --- new "dmg" column (T), only on rows that carry an "irl price" (Q) value ---
$ws.Range("T2").Value = "dmg"
$rowsWithQ = @(3,5,6,7,9,10,12,13,14,15,16,17,18,20,21,23,24,26,27,28,30,31)
foreach ($r in $rowsWithQ) {
    $ws.Range("T$r").Formula = "=(Q$r-5)*0.09/11"
}

# --- restore selection to match the edited cell ---
[void]$ws.Range("T3").Select()
